$d = $word.ActiveDocument

foreach ($t in $d.Tables) {
    $t.Borders.Enable = $true
    $t.Borders.OutsideLineStyle = 1
    $t.Borders.InsideLineStyle = 1
    $t.Borders.OutsideLineWidth = 2
    $t.Borders.InsideLineWidth = 2
    $t.Borders.OutsideColor = 0
    $t.Borders.InsideColor = 0
}

Write-Output "Tables updated: $($d.Tables.Count)"
